$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recreate the frozen header pane (row 1 frozen) -------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Apply the "data row" font styling (Calibri 10, automatic/theme text) --
# to every cell that will hold data in rows 2-7, skipping column B (which the
# source data never touches) so it is not given a style.
$rngData = $ws.Range("A2:A7,C2:H7")
foreach ($area in $rngData.Areas) {
  $area.Font.Name = "Calibri"
  $area.Font.Size = 10
  $area.Font.ThemeColor = 1
}

# --- Fill in the new descriptive rows ---------------------------------------
$ws.Range("A2").Value = "MCH160-1"
$ws.Range("C2").Value = " NEWSLETTER 1985-1992"
$ws.Range("D2").Value = "1985-1992"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 22B | GRAP COUNT NUMER: NONE"

$ws.Range("A3").Value = "MCH160-2"
$ws.Range("C3").Value = "ANNUAL REPORTS, UN REPORTS, SHIREBU"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 22B | GRAP COUNT NUMER: NONE"

$ws.Range("A4").Value = "MCH160-3"
$ws.Range("C4").Value = "VARIOUS PUBLICATIONS, UN REPORTS"
$ws.Range("E4").Value = "Series"
$ws.Range("F4").Value = "1 Box"
$ws.Range("G4").Value = "LOCATION: 22B | GRAP COUNT NUMER: NONE"

$ws.Range("A5").Value = "MCH160-4"
$ws.Range("C5").Value = "SHIREBU'S REPORTS"
$ws.Range("E5").Value = "Series"
$ws.Range("F5").Value = "1 Box"
$ws.Range("G5").Value = "LOCATION: 22B | GRAP COUNT NUMER: NONE"

$ws.Range("A6").Value = "MCH160-5"
$ws.Range("C6").Value = "SHIREBU'S REPORTS, NEWSLETTERS"
$ws.Range("E6").Value = "Series"
$ws.Range("F6").Value = "1 Box"
$ws.Range("G6").Value = "LOCATION: 22B | GRAP COUNT NUMER: NONE"

$ws.Range("A7").Value = "MCH160-6"
$ws.Range("C7").Value = "VARIOUS PUBLICATIONS  "
$ws.Range("E7").Value = "Series"
$ws.Range("F7").Value = "1 Box"
$ws.Range("G7").Value = "LOCATION: 22B | GRAP COUNT NUMER: NONE"

# --- Final UI selection (matches the saved view in the workbook) -----------
$ws.Range("C13").Select()
